$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 644, shifting existing rows 644:685 down to 645:686
$ws.Rows.Item(644).Insert()

# Populate the newly inserted row 644 with the new data.
# The date-shaped text must stay a literal text string (matching the
# original inlineStr cells), so prefix with an apostrophe to stop Excel
# from auto-converting it into a date serial value.
$ws.Cells.Item(644, 1).Value = "'2026/01/13"
$ws.Cells.Item(644, 2).Value = "火"
$ws.Cells.Item(644, 3).Value = 23
$ws.Cells.Item(644, 4).Value = 201
